# Apply the "functional model need to fix variables" fix to the
# CustomerVehicleTripMap workbook. Each sheet (y1..y8) is a 10x6 (A1:F10)
# 0/1 "selection" matrix. This patch flips a handful of individual cells
# across several sheets.

$wb = $excel.ActiveWorkbook

# --- y1 ---
$ws = $wb.Worksheets.Item("y1")
$ws.Range("A7").Value = 1

# --- y3 ---
$ws = $wb.Worksheets.Item("y3")
$ws.Range("C4").Value = 1

# --- y4 ---
$ws = $wb.Worksheets.Item("y4")
$ws.Range("E8").Value = 1

# --- y5 ---
$ws = $wb.Worksheets.Item("y5")
$ws.Range("A2").Value = 1
$ws.Range("C6").Value = 0.9999995974647448

# --- y6 ---
$ws = $wb.Worksheets.Item("y6")
$ws.Range("E6").Value = 0.0000004014241441439463

# --- y8 ---
$ws = $wb.Worksheets.Item("y8")
$ws.Range("C2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("C3").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("B5").Value = 0
$ws.Range("F5").Value = 1
$ws.Range("E6").Value = 0
$ws.Range("A7").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("B9").Value = 1
$ws.Range("D9").Value = 0
